# Appends new XRPUSDT 4h OHLCV rows (1086-1120) to Sheet1, matching the
# commit that extends the dataset from row 1085 through row 1120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1086
$endRow = 1120

# New rows: Timestamp(serial), Open, High, Low, Close, Volume
$data = @(
    @(45534.66666666666, 0.5485, 0.5629, 0.5463, 0.5618, 47431744),
    @(45534.83333333334, 0.5617, 0.5629999999999999, 0.5603, 0.5615, 2608141),
    @(45535, 0.5667, 0.5682, 0.5653, 0.5663, 4516773),
    @(45535.16666666666, 0.5648, 0.5687, 0.5645, 0.5665, 8198539),
    @(45535.33333333334, 0.5665, 0.5721000000000001, 0.5625, 0.5692, 17674381),
    @(45535.5, 0.5693, 0.57, 0.5644, 0.5664, 12896839),
    @(45535.66666666666, 0.5664, 0.5681, 0.5624, 0.5678, 11118572),
    @(45535.83333333334, 0.5678, 0.5681, 0.5661, 0.5671, 2299466),
    @(45536, 0.5662, 0.5667, 0.5590000000000001, 0.5608, 11925034),
    @(45536.16666666666, 0.5607, 0.5607, 0.553, 0.5591, 16206420),
    @(45536.33333333334, 0.5590000000000001, 0.5603, 0.5580000000000001, 0.5590000000000001, 4666214),
    @(45536.5, 0.5581, 0.5589, 0.5515, 0.5565, 22094783),
    @(45536.66666666666, 0.5566, 0.5583, 0.5558, 0.5558, 2045733),
    @(45536.83333333334, 0.5595, 0.5619, 0.5553, 0.5579, 5416155),
    @(45537, 0.5471, 0.5508, 0.5463, 0.5501, 30593301),
    @(45537.16666666666, 0.5501, 0.5511, 0.5443, 0.5463, 34732817),
    @(45537.33333333334, 0.5463, 0.5574, 0.5455, 0.5562, 24537945),
    @(45537.5, 0.5548999999999999, 0.5561, 0.554, 0.5547, 5216991),
    @(45537.66666666666, 0.5582, 0.5679999999999999, 0.5555, 0.5679, 22420724),
    @(45537.83333333334, 0.5678, 0.5679, 0.5656, 0.5671, 2857569),
    @(45538, 0.5677, 0.5724, 0.5656, 0.5667, 24133935),
    @(45538.16666666666, 0.5666, 0.5695, 0.5649, 0.5692, 31112117),
    @(45538.33333333334, 0.5691000000000001, 0.5724, 0.5658, 0.5701000000000001, 30808267),
    @(45538.5, 0.57, 0.5723, 0.5571, 0.5593, 50553028),
    @(45538.66666666666, 0.5592, 0.5659, 0.5586, 0.5659, 25331286),
    @(45538.83333333334, 0.5659, 0.5677, 0.5565, 0.5572, 17842850),
    @(45539, 0.5572, 0.5604, 0.5326, 0.553, 95792961),
    @(45539.16666666666, 0.553, 0.5556, 0.5497, 0.5545, 33550540),
    @(45539.33333333334, 0.5546, 0.5570000000000001, 0.549, 0.5505, 43266166),
    @(45539.5, 0.5505, 0.5629999999999999, 0.5472, 0.5621, 58081896),
    @(45539.66666666666, 0.5621, 0.5648, 0.5564, 0.5567, 22582496),
    @(45539.83333333334, 0.5578, 0.5580000000000001, 0.5574, 0.5578, 118795),
    @(45540, 0.5579, 0.5598, 0.5517, 0.5528, 29278558),
    @(45540.16666666666, 0.5528999999999999, 0.5547, 0.5499000000000001, 0.5535, 35982005),
    @(45540.33333333334, 0.5535, 0.5547, 0.5498, 0.55, 6988366)
)

# Reuse the existing date/time style (border + bold + center + custom
# numFmt) from the last existing timestamp cell (A1085) for the new
# timestamp cells in column A, so no new style entries are created.
$ws.Range("A1085").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowValues[0]
    $ws.Cells.Item($r, 2).Value = $rowValues[1]
    $ws.Cells.Item($r, 3).Value = $rowValues[2]
    $ws.Cells.Item($r, 4).Value = $rowValues[3]
    $ws.Cells.Item($r, 5).Value = $rowValues[4]
    $ws.Cells.Item($r, 6).Value = $rowValues[5]
}

Write-Output "Appended $($data.Count) rows ($startRow to $endRow)"
